$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update values in rows 4-58 (unaffected by the row deletion below) ---
$ws.Range("D4").Value = 2448
$ws.Range("E4").Value = 10603

$ws.Range("D5").Value = "1015 (41.5)"
$ws.Range("E5").Value = "3965 (37.4)"

$ws.Range("D6").Value = "596 (24.3)"
$ws.Range("E6").Value = "2561 (24.2)"

$ws.Range("D7").Value = "507 (20.7)"
$ws.Range("E7").Value = "2420 (22.8)"

$ws.Range("D8").Value = "330 (13.5)"
$ws.Range("E8").Value = "1657 (15.6)"

$ws.Range("D9").Value = "1170 (47.8)"
$ws.Range("E9").Value = "4450 (42.0)"

$ws.Range("D10").Value = "932 (38.1)"
$ws.Range("E10").Value = "5370 (50.6)"

$ws.Range("D11").Value = "342 (14.0)"
$ws.Range("E11").Value = "560 (5.3)"

$ws.Range("D12").Value = "1174 (48.0)"
$ws.Range("E12").Value = "4673 (44.1)"

$ws.Range("D13").Value = "770 (31.5)"
$ws.Range("E13").Value = "490 (4.6)"

$ws.Range("D14").Value = "196 (8.0)"
$ws.Range("E14").Value = "1545 (14.6)"

$ws.Range("D15").Value = "805 (32.9)"
$ws.Range("E15").Value = "4481 (42.3)"

$ws.Range("D16").Value = "434 (17.7)"
$ws.Range("E16").Value = "1762 (16.6)"

$ws.Range("D17").Value = "650 (26.6)"
$ws.Range("E17").Value = "2667 (25.2)"

$ws.Range("D18").Value = "1207 (49.3)"
$ws.Range("E18").Value = "5398 (50.9)"

$ws.Range("D19").Value = "227 (9.3)"
$ws.Range("E19").Value = "494 (4.7)"

$ws.Range("D20").Value = "999 (40.8)"
$ws.Range("E20").Value = "5074 (47.9)"

$ws.Range("D21").Value = "1392 (56.9)"
$ws.Range("E21").Value = "6318 (59.6)"

$ws.Range("D22").Value = "371 (15.2)"
$ws.Range("E22").Value = "1027 (9.7)"

$ws.Range("D23").Value = "1181 (48.2)"
$ws.Range("E23").Value = "5932 (55.9)"

$ws.Range("D24").Value = "1502 (61.4)"
$ws.Range("E24").Value = "6318 (59.6)"

$ws.Range("D25").Value = "256 (10.5)"
$ws.Range("E25").Value = "1296 (12.2)"

$ws.Range("D26").Value = "2385 (97.4)"
$ws.Range("E26").Value = "10357 (97.7)"

$ws.Range("D27").Value = "1704 (69.6)"
$ws.Range("E27").Value = "6977 (65.8)"

$ws.Range("D28").Value = "915 (37.4)"
$ws.Range("E28").Value = "3851 (36.3)"

$ws.Range("D29").Value = "562 (23.0)"
$ws.Range("E29").Value = "2623 (24.7)"

$ws.Range("D30").Value = "34 (1.4)"
$ws.Range("E30").Value = "152 (1.4)"

$ws.Range("D31").Value = "754 (30.8)"
$ws.Range("E31").Value = "4137 (39.0)"

$ws.Range("D32").Value = "3 (0.1)"

$ws.Range("D33").Value = "19 (0.8)"
$ws.Range("E33").Value = "68 (0.6)"

$ws.Range("D34").Value = "161 (6.6)"
$ws.Range("E34").Value = "532 (5.0)"

$ws.Range("D35").Value = "62 (2.5)"
$ws.Range("E35").Value = "185 (1.7)"

$ws.Range("D36").Value = "264 (10.8)"
$ws.Range("E36").Value = "461 (4.3)"

$ws.Range("D37").Value = "1939 (79.2)"
$ws.Range("E37").Value = "9355 (88.2)"

$ws.Range("D38").Value = "75 (3.1)"
$ws.Range("E38").Value = "300 (2.8)"

$ws.Range("D39").Value = "986 (40.3)"
$ws.Range("E39").Value = "2962 (27.9)"

$ws.Range("D40").Value = "1387 (56.7)"
$ws.Range("E40").Value = "7341 (69.2)"

$ws.Range("D41").Value = "117 (4.8)"
$ws.Range("E41").Value = "472 (4.5)"

$ws.Range("D42").Value = "139 (5.7)"
$ws.Range("E42").Value = "507 (4.8)"

$ws.Range("D43").Value = "15 (0.6)"
$ws.Range("E43").Value = "84 (0.8)"

$ws.Range("D44").Value = "2 (0.1)"
$ws.Range("E44").Value = "20 (0.2)"

$ws.Range("D45").Value = "4 (0.2)"

$ws.Range("D46").Value = "64 [52,75]"
$ws.Range("E46").Value = "68 [58,78]"

$ws.Range("C47").Value = 10855
$ws.Range("D47").Value = "5.23 [3.23,8.99]"
$ws.Range("E47").Value = "5.38 [3.29,9.38]"

$ws.Range("C48").Value = 2196
$ws.Range("D48").Value = "4.17 [2.83,7.71]"
$ws.Range("E48").Value = "4.08 [2.79,7.17]"

$ws.Range("C49").Value = 10855
$ws.Range("D49").Value = "10.00 [6.00,19.00]"
$ws.Range("E49").Value = "11.00 [6.00,18.75]"

$ws.Range("C50").Value = 2196

$ws.Range("D52").Value = "6 [4,8]"
$ws.Range("E52").Value = "5 [3,8]"

$ws.Range("C53").Value = 4756

$ws.Range("C54").Value = 28

$ws.Range("C55").Value = 5163

$ws.Range("C56").Value = 18
$ws.Range("D56").Value = "1 [1,3]"

$ws.Range("C57").Value = 22

# --- Remove the 7 rows for *_max_0_24h (SOFA sub-scores 0-24h), rows 59-65 ---
# This shifts old rows 66-91 up to become new rows 59-84,
# and also updates the sheet dimension from A1:E91 to A1:E84.
$ws.Range("A59:E65").EntireRow.Delete()

# --- Update values in rows 59-84 (post-shift; these were old rows 66-91) ---
$ws.Range("C59").Value = 4731
$ws.Range("D59").Value = "729 [246,1498]"
$ws.Range("E59").Value = "650 [210,1463]"

$ws.Range("C60").Value = 309
$ws.Range("D60").Value = "2648 [994,5920]"
$ws.Range("E60").Value = "2369 [817,5431]"

$ws.Range("C61").Value = 309
$ws.Range("D61").Value = "506.2 [226.7,963.1]"
$ws.Range("E61").Value = "477.0 [204.8,905.2]"

$ws.Range("C62").Value = 7672
$ws.Range("D62").Value = "50 [40,60]"
$ws.Range("E62").Value = "50 [42,65]"

$ws.Range("C63").Value = 5341
$ws.Range("D63").Value = "39.0 [17.0,86.2]"
$ws.Range("E63").Value = "32.0 [15.0,76.0]"

$ws.Range("C64").Value = 5341
$ws.Range("D64").Value = "0.30 [0.14,0.50]"
$ws.Range("E64").Value = "0.26 [0.13,0.45]"

$ws.Range("C65").Value = 5341
$ws.Range("D65").Value = "2.0 [1.0,8.0]"
$ws.Range("E65").Value = "3.0 [1.0,8.0]"

$ws.Range("C66").Value = 11653
$ws.Range("D66").Value = "16.0 [4.0,41.0]"
$ws.Range("E66").Value = "25.0 [5.5,63.0]"

$ws.Range("C67").Value = 5827
$ws.Range("D67").Value = "3.0 [1.0,12.0]"
$ws.Range("E67").Value = "3.0 [1.0,10.0]"

$ws.Range("C68").Value = 5827
$ws.Range("D68").Value = "29.0 [10.0,65.0]"
$ws.Range("E68").Value = "29.0 [11.0,62.0]"

$ws.Range("C69").Value = 5827
$ws.Range("D69").Value = "0.23 [0.08,0.48]"
$ws.Range("E69").Value = "0.25 [0.09,0.47]"

$ws.Range("C70").Value = 24
$ws.Range("D70").Value = "19.4 [17.0,22.6]"
$ws.Range("E70").Value = "19.1 [16.8,22.0]"

$ws.Range("C71").Value = 18
$ws.Range("D71").Value = "77.5 [71.1,85.4]"
$ws.Range("E71").Value = "74.7 [69.4,81.6]"

$ws.Range("C72").Value = 693
$ws.Range("D72").Value = "36.9 [36.6,37.3]"
$ws.Range("E72").Value = "36.9 [36.6,37.2]"

$ws.Range("C73").Value = 19
$ws.Range("D73").Value = "98.0 [96.3,99.2]"
$ws.Range("E73").Value = "97.2 [95.8,98.6]"

$ws.Range("C74").Value = 18
$ws.Range("D74").Value = "87.5 [76.4,100.6]"
$ws.Range("E74").Value = "85.8 [75.8,97.9]"

$ws.Range("C75").Value = 4014
$ws.Range("D75").Value = "88.5 [69.0,124.0]"
$ws.Range("E75").Value = "89.0 [72.0,118.0]"

$ws.Range("C76").Value = 4014
$ws.Range("D76").Value = "44.0 [37.0,52.0]"
$ws.Range("E76").Value = "46.0 [39.0,53.0]"

$ws.Range("C77").Value = 2267
$ws.Range("D77").Value = "7.3 [7.2,7.4]"
$ws.Range("E77").Value = "7.3 [7.2,7.4]"

$ws.Range("C78").Value = 67
$ws.Range("D78").Value = "154.0 [122.0,218.0]"
$ws.Range("E78").Value = "147.0 [120.0,194.0]"

$ws.Range("C79").Value = 27
$ws.Range("D79").Value = "137.0 [134.0,140.0]"
$ws.Range("E79").Value = "137.0 [134.0,140.0]"

$ws.Range("C80").Value = 34
$ws.Range("D80").Value = "4.5 [4.1,5.2]"
$ws.Range("E80").Value = "4.5 [4.1,5.0]"

$ws.Range("C81").Value = 12719
$ws.Range("D81").Value = "19.0 [13.5,30.9]"
$ws.Range("E81").Value = "22.0 [12.5,33.9]"

$ws.Range("C82").Value = 1753
$ws.Range("D82").Value = "9.8 [8.2,11.4]"
$ws.Range("E82").Value = "10.0 [8.6,11.6]"

$ws.Range("C83").Value = 8424
$ws.Range("D83").Value = "226.0 [153.0,364.0]"
$ws.Range("E83").Value = "228.5 [165.0,344.0]"

$ws.Range("C84").Value = 788
$ws.Range("D84").Value = "1.4 [1.2,1.7]"
$ws.Range("E84").Value = "1.4 [1.2,1.7]"
